$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet's conversion summary text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.33 = 29708.57 pesos`n✅ 29708.57 pesos = 7.3 = 961.92 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 136.5
$wsTasas.Range("O10").Value = 4055.22
$wsTasas.Range("N12").Value = 4070
$wsTasas.Range("O12").Value = 131.78
